$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("G28").Value = 1.7
$ws.Range("I28").Value = 5.75
$ws.Range("J28").Value = 2.4
$ws.Range("Q28").Value = 2.5
$ws.Range("R28").Value = 1.5
$ws.Range("U28").Value = 5
$ws.Range("V28").Value = 1.17
$ws.Range("Y28").Value = 2.38
$ws.Range("Z28").Value = 1.53
$ws.Range("AA28").Value = 5
$ws.Range("AI28").Value = 23
$ws.Range("AM28").Value = 26
$ws.Range("AN28").Value = 19
$ws.Range("AR28").Value = 1.95
$ws.Range("AS28").Value = 1.9

# Row 29
$ws.Range("G29").Value = 2.63
$ws.Range("I29").Value = 2.9
$ws.Range("K29").Value = 1.83
$ws.Range("O29").Value = 1.57
$ws.Range("P29").Value = 2.25
$ws.Range("S29").Value = 4.8
$ws.Range("T29").Value = 1.19
$ws.Range("U29").Value = 6
$ws.Range("V29").Value = 1.13
$ws.Range("AA29").Value = 6
$ws.Range("AD29").Value = 26
$ws.Range("AE29").Value = 29
$ws.Range("AG29").Value = 6
$ws.Range("AI29").Value = 21
$ws.Range("AN29").Value = 12
$ws.Range("AR29").Value = 2.05
$ws.Range("AS29").Value = 1.8

# Row 41
$ws.Range("M41").Value = 1.13
$ws.Range("N41").Value = 6
$ws.Range("O41").Value = 1.53
$ws.Range("P41").Value = 2.38
$ws.Range("Q41").Value = 2.7
$ws.Range("R41").Value = 1.44
$ws.Range("U41").Value = 5.5
$ws.Range("V41").Value = 1.14
$ws.Range("W41").Value = 1.62
$ws.Range("X41").Value = 2.2
$ws.Range("Y41").Value = 2.2
$ws.Range("Z41").Value = 1.62
$ws.Range("AE41").Value = 34
$ws.Range("AI41").Value = 19
$ws.Range("AL41").Value = 6
$ws.Range("AP41").Value = 26
$ws.Range("AR41").Value = 2.03
$ws.Range("AS41").Value = 1.83

# Row 77
$ws.Range("G77").Value = 3
$ws.Range("I77").Value = 2.4
$ws.Range("J77").Value = 3.75
$ws.Range("K77").Value = 2
$ws.Range("L77").Value = 3.2
$ws.Range("Q77").Value = 2.1
$ws.Range("R77").Value = 1.7
$ws.Range("U77").Value = 3.75
$ws.Range("V77").Value = 1.25
$ws.Range("W77").Value = 1.5
$ws.Range("X77").Value = 2.5
$ws.Range("Y77").Value = 1.95
$ws.Range("Z77").Value = 1.8
$ws.Range("AD77").Value = 34
$ws.Range("AE77").Value = 26
$ws.Range("AF77").Value = 41
$ws.Range("AG77").Value = 8
$ws.Range("AI77").Value = 17
$ws.Range("AK77").Value = 351
$ws.Range("AL77").Value = 7
$ws.Range("AM77").Value = 11
$ws.Range("AO77").Value = 23

# Row 80
$ws.Range("G80").Value = 3.55
$ws.Range("H80").Value = 3.35
$ws.Range("I80").Value = 1.95
$ws.Range("J80").Value = 3.95
$ws.Range("K80").Value = 2.1
$ws.Range("L80").Value = 2.57
$ws.Range("Q80").Value = 1.8
$ws.Range("U80").Value = 2.85
$ws.Range("V80").Value = 1.32
$ws.Range("W80").Value = 1.38
$ws.Range("X80").Value = 2.57
$ws.Range("AA80").Value = 11.25
$ws.Range("AB80").Value = 20
$ws.Range("AC80").Value = 12
$ws.Range("AD80").Value = 50
$ws.Range("AE80").Value = 30
$ws.Range("AF80").Value = 35
$ws.Range("AG80").Value = 10.25
$ws.Range("AH80").Value = 6.5
$ws.Range("AL80").Value = 7.4
$ws.Range("AM80").Value = 9.5
$ws.Range("AN80").Value = 8.5
$ws.Range("AO80").Value = 17.5
$ws.Range("AQ80").Value = 26
